$p = $ppt.ActivePresentation
try {
  $cp = $p.CustomXMLParts
  Write-Output "count=$($cp.Count)"
} catch { Write-Output "ERR $_" }
